# Robo_Extrator/ID/Identidades/Registro.xlsx
# "Iniciando desafio robô registro identidade"
#
# The diff duplicates the single data row in the "PESSOAL" sheet: a new
# row 3 is inserted holding the same identity record (A: id number,
# B: name, C: birth date) that already lives in row 2.
#
# We insert the new row *below* the existing data (at row 3) so the blank
# row inherits row 2's formatting (date style on column C) automatically,
# then fill it with the same values as row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 3; it inherits formatting from row 2 above it.
$ws.Rows.Item(3).Insert()

# A3: same numeric ID as A2.
$ws.Cells.Item(3, 1).Value = $ws.Cells.Item(2, 1).Value2

# B3: same name text as B2 (copy/paste keeps it an exact duplicate string
# and avoids re-triggering row autofit from a literal multi-line value).
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(3, 2).PasteSpecial(-4163)

# C3: same birth-date serial value as C2.
$ws.Cells.Item(3, 3).Value = $ws.Cells.Item(2, 3).Value2
